$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(42602.010092592594, "Bag", 5860, 6487, 815, 130, 59, 67, 30, 0, 0, 0, 0),
    @(42602.481874999998, "Bag", 5803, 9336, 1178, 193, 101, 65, 34, 1, 0, 100, 0)
)

$r = 17
foreach ($row in $newRows) {
    $ws.Cells.Item(16, 1).Copy($ws.Cells.Item($r, 1))
    $ws.Cells.Item($r, 1).Value = $row[0]
    for ($c = 2; $c -le 13; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}
